# Adapt column header formatting to respective input file names:
#   *_old  -> *_FV2410
#   *_new  -> *_FV2504
# Wrap the data range in an Excel Table and freeze the header row.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$fv2410 = "_FV2410"
$fv2504 = "_FV2504"

$baseNames = @("Segmentname", "Segmentgruppe", "Segment", "Datenelement", "Segment ID", "Code", "Qualifier", "Beschreibung", "Bedingungsausdruck", "Bedingung")

# Columns A-J (1-10): "<name>_old" -> "<name>_FV2410"
for ($i = 0; $i -lt $baseNames.Count; $i++) {
    $col = $i + 1
    $ws.Cells.Item(1, $col).Value2 = $baseNames[$i] + $fv2410
}

# Column K (11): "diff" stays unchanged

# Columns L-U (12-21): "<name>_new" -> "<name>_FV2504"
for ($i = 0; $i -lt $baseNames.Count; $i++) {
    $col = $i + 12
    $ws.Cells.Item(1, $col).Value2 = $baseNames[$i] + $fv2504
}

# Turn the A1:U71 range into a native Excel table.
$tableRange = $ws.Range("A1:U71")
$tbl = $ws.ListObjects.Add(1, $tableRange, $null, 1)
$tbl.Name = "Table1"

# Freeze the header row (split after row 1, keep top-left at A2).
$null = $ws.Range("A2").Select()
$excel.ActiveWindow.FreezePanes = $true
